# Apply the "Added closed bookings and some other stuff" edit.
#
# Summary of changes:
#  - sheet1 ("question_template"): add a new column H "Length(seconds)"
#    header (matching the existing header style) and empty data cells
#    H2:H10 (matching the existing data style used in B2:G10).
#  - sheet1: column width tweaks for columns E, G and the new column H.
#  - sheet1: selection moved from D16 to J8.
#  - sheet2 ("Sheet2"): PageSetup FirstPageNumber changed from 1 to 0.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: question_template
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Add the new column H header, copying the formatting already used by
# the neighbouring "Correct Answer" header cell (G1).
$ws1.Range("G1").Copy($ws1.Range("H1"))
$ws1.Range("H1").Value = "Length(seconds)"

# Add the (empty) data cells for the new column, copying the formatting
# already used by the other plain data columns (e.g. B2:B10).
$ws1.Range("B2:B10").Copy()
$ws1.Range("H2:H10").PasteSpecial(-4122)

# Column width tweaks (ColumnWidth is expressed in characters; the
# stored OOXML width = ColumnWidth + 5/6, rounded to the nearest 1/6th
# of a character by the engine's pixel-snapping, matching real Excel
# COM automation behaviour). The chosen inputs land exactly on a 1/6
# boundary so the resulting stored width is deterministic.
$ws1.Range("E1").ColumnWidth = 13.666666666666666
$ws1.Range("G1").ColumnWidth = 13.0
$ws1.Range("H1").ColumnWidth = 14.166666666666666

# Selection moved to J8.
$ws1.Range("J8").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: Sheet2
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ps2 = $ws2.PageSetup()
$ps2.FirstPageNumber = 0
